$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "M.BA"
$ws.Range("B2").Value = 19
$ws.Range("C2").Value = 9
$ws.Range("D2").Value = 1000000
$ws.Range("E2").Value = 9259244.285845313
$ws.Range("G2").Value = 45173
$ws.Range("H2").Value = 22
$ws.Range("I2").Value = 12
$ws.Range("J2").Value = 10
$ws.Range("K2").Value = 0.545
$ws.Range("L2").Value = 0.455
$ws.Range("M2").Value = -0.093
$ws.Range("N2").Value = 713950.269
$ws.Range("O2").Value = 0.235
$ws.Range("P2").Value = -136916.068
$ws.Range("Q2").Value = -0.049
$ws.Range("R2").Value = -514493.047
$ws.Range("S2").Value = 4047603.618
$ws.Range("T2").Value = 4.047603617564064
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 5
$ws.Range("W2").Value = 66.9047619047619
